# run_metrics.xlsx update: results re-synced from the DFKI cluster run.
# Rows are re-ordered/re-numbered (2..45) to reflect the freshly synced runs,
# and the sheet grows from 33 to 44 data rows (dimension A1:C34 -> A1:C45).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "best accuracy" highlight (yellow) moves from the old row 27
# (U_Net_LFCC_32_len30S) to the new row 7 (U_Net_LFCC_32_len5S).
# Copy/paste the format so the existing fill/style is reused, then clear
# the donor cell's formatting.
$ws.Range("B27").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("B27").ClearFormats()
$excel.CutCopyMode = $false

# Refresh every data row (A: run name, B: best accuracy, C: least loss).
# A handful of runs legitimately report 0 accuracy / infinite loss; those
# are written as literal text "inf" in column C, matching the source data.
$ws.Range("A2").Value = 'R2AttU_Net_lfcc-delta_32_len5S'
$ws.Range("B2").Value = 69.56521739130434
$ws.Range("C2").Value = 0.557578980922699
$ws.Range("A3").Value = 'U_Net_delta_32_len30S'
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = "inf"
$ws.Range("A4").Value = 'U_Net_delta_32_len30S'
$ws.Range("B4").Value = 93.47826086956522
$ws.Range("C4").Value = 0.8856147925059
$ws.Range("A5").Value = 'U_Net_delta_32_len30S'
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = "inf"
$ws.Range("A6").Value = 'U_Net_LFCC_80_len5S'
$ws.Range("B6").Value = 91.30434782608695
$ws.Range("C6").Value = 0.8510100543498993
$ws.Range("A7").Value = 'U_Net_LFCC_32_len5S'
$ws.Range("B7").Value = 95.65217391304348
$ws.Range("C7").Value = 0.8168511788050333
$ws.Range("A8").Value = 'AttentionUNet_lfcc-delta-delta_32_len5S'
$ws.Range("B8").Value = 76.08695652173913
$ws.Range("C8").Value = 0.5528329908847809
$ws.Range("A9").Value = 'AttentionUNet_MFCC_32_len5S'
$ws.Range("B9").Value = 67.39130434782609
$ws.Range("C9").Value = 0.5642368793487549
$ws.Range("A10").Value = 'U_Net_MFCC_32_len5S'
$ws.Range("B10").Value = 86.95652173913044
$ws.Range("C10").Value = 0.9281513889630636
$ws.Range("A11").Value = 'AttentionUNet_delta_32_len5S'
$ws.Range("B11").Value = 82.60869565217391
$ws.Range("C11").Value = 0.5391703248023987
$ws.Range("A12").Value = 'R2U_Net_LFCC_32_len5S'
$ws.Range("B12").Value = 73.91304347826087
$ws.Range("C12").Value = 0.6428556889295578
$ws.Range("A13").Value = 'U_Net_lfcc-delta-delta_32_len5S'
$ws.Range("B13").Value = 78.26086956521739
$ws.Range("C13").Value = 0.9822227358818054
$ws.Range("A14").Value = 'R2U_Net_delta_32_len5S'
$ws.Range("B14").Value = 73.91304347826087
$ws.Range("C14").Value = 0.5048011690378189
$ws.Range("A15").Value = 'R2AttU_Net_MFCC_32_len5S'
$ws.Range("B15").Value = 76.08695652173913
$ws.Range("C15").Value = 0.5616060495376587
$ws.Range("A16").Value = 'U_Net_delta_32_len5S'
$ws.Range("B16").Value = 71.73913043478261
$ws.Range("C16").Value = 1.076469659805298
$ws.Range("A17").Value = 'R2AttU_Net_delta_32_len5S'
$ws.Range("B17").Value = 69.56521739130434
$ws.Range("C17").Value = 0.5860709249973297
$ws.Range("A18").Value = 'U_Net_delta_80_len5S'
$ws.Range("B18").Value = 63.04347826086956
$ws.Range("C18").Value = 1.033747345209122
$ws.Range("A19").Value = 'U_Net_MFCC_32_len30S'
$ws.Range("B19").Value = 89.1304347826087
$ws.Range("C19").Value = 0.8510035673777262
$ws.Range("A20").Value = 'U_Net_MFCC_32_len30S'
$ws.Range("B20").Value = 86.95652173913044
$ws.Range("C20").Value = 0.8733596205711365
$ws.Range("A21").Value = 'U_Net_MFCC_32_len30S'
$ws.Range("B21").Value = 0
$ws.Range("C21").Value = "inf"
$ws.Range("A22").Value = 'AttentionUNet_LFCC_32_len5S'
$ws.Range("B22").Value = 82.60869565217391
$ws.Range("C22").Value = 0.4653773456811905
$ws.Range("A23").Value = 'R2AttU_Net_lfcc-delta-delta_32_len5S'
$ws.Range("B23").Value = 54.34782608695652
$ws.Range("C23").Value = 0.6629349291324615
$ws.Range("A24").Value = 'U_Net_lfcc-delta-delta_80_len5S'
$ws.Range("B24").Value = 89.1304347826087
$ws.Range("C24").Value = 0.9213963150978088
$ws.Range("A25").Value = 'R2U_Net_MFCC_32_len5S'
$ws.Range("B25").Value = 50
$ws.Range("C25").Value = 0.6707836836576462
$ws.Range("A26").Value = 'U_Net_delta-delta_32_len30S'
$ws.Range("B26").Value = 78.26086956521739
$ws.Range("C26").Value = 0.9564621647198995
$ws.Range("A27").Value = 'U_Net_delta-delta_32_len30S'
$ws.Range("B27").Value = 0
$ws.Range("C27").Value = "inf"
$ws.Range("A28").Value = 'U_Net_lfcc-delta-delta_32_len30S'
$ws.Range("B28").Value = 0
$ws.Range("C28").Value = "inf"
$ws.Range("A29").Value = 'U_Net_lfcc-delta-delta_32_len30S'
$ws.Range("B29").Value = 63.04347826086956
$ws.Range("C29").Value = 1.059326450030009
$ws.Range("A30").Value = 'U_Net_delta-delta_32_len5S'
$ws.Range("B30").Value = 63.04347826086956
$ws.Range("C30").Value = 1.100150108337402
$ws.Range("A31").Value = 'U_Net_lfcc-delta_32_len30S'
$ws.Range("B31").Value = 0
$ws.Range("C31").Value = "inf"
$ws.Range("A32").Value = 'U_Net_lfcc-delta_32_len30S'
$ws.Range("B32").Value = 84.78260869565217
$ws.Range("C32").Value = 0.9093547463417053
$ws.Range("A33").Value = 'U_Net_LFCC_32_len30S'
$ws.Range("B33").Value = 86.95652173913044
$ws.Range("C33").Value = 0.8723922967910767
$ws.Range("A34").Value = 'U_Net_LFCC_32_len30S'
$ws.Range("B34").Value = 0
$ws.Range("C34").Value = "inf"
$ws.Range("A35").Value = 'R2U_Net_lfcc-delta_32_len5S'
$ws.Range("B35").Value = 78.26086956521739
$ws.Range("C35").Value = 0.5619004666805267
$ws.Range("A36").Value = 'R2AttU_Net_delta-delta_32_len5S'
$ws.Range("B36").Value = 58.69565217391305
$ws.Range("C36").Value = 0.665488988161087
$ws.Range("A37").Value = 'R2U_Net_lfcc-delta-delta_32_len5S'
$ws.Range("B37").Value = 67.39130434782609
$ws.Range("C37").Value = 0.5734604299068451
$ws.Range("A38").Value = 'R2AttU_Net_LFCC_32_len5S'
$ws.Range("B38").Value = 73.91304347826087
$ws.Range("C38").Value = 0.533536896109581
$ws.Range("A39").Value = 'AttentionUNet_lfcc-delta_32_len5S'
$ws.Range("B39").Value = 78.26086956521739
$ws.Range("C39").Value = 0.5578003972768784
$ws.Range("A40").Value = 'U_Net_lfcc-delta_80_len5S'
$ws.Range("B40").Value = 93.47826086956522
$ws.Range("C40").Value = 0.8924888968467712
$ws.Range("A41").Value = 'R2U_Net_delta-delta_32_len5S'
$ws.Range("B41").Value = 73.91304347826087
$ws.Range("C41").Value = 0.621788889169693
$ws.Range("A42").Value = 'U_Net_delta-delta_80_len5S'
$ws.Range("B42").Value = 50
$ws.Range("C42").Value = 1.081704556941986
$ws.Range("A43").Value = 'AttentionUNet_delta-delta_32_len5S'
$ws.Range("B43").Value = 60.8695652173913
$ws.Range("C43").Value = 0.5899332761764526
$ws.Range("A44").Value = 'U_Net_MFCC_80_len5S'
$ws.Range("B44").Value = 78.26086956521739
$ws.Range("C44").Value = 0.9191123247146606
$ws.Range("A45").Value = 'U_Net_lfcc-delta_32_len5S'
$ws.Range("B45").Value = 67.39130434782609
$ws.Range("C45").Value = 1.047885318597158

# Highlight the new "least loss" best result (C22, AttentionUNet_LFCC_32_len5S)
# with a light-green fill (#90EE90).
$ws.Range("C22").Interior.Color = 9498256
